$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Range("D$r").Value = 0.0582
    $ws.Range("G$r").Value = 0.2593175853018373
    $ws.Range("H$r").Value = 0.2593175853018373
    $ws.Range("I$r").Value = -0.07769028871391076
    $ws.Range("J$r").Value = -0.07769028871391076
    $ws.Range("K$r").Value = -5.12
    $ws.Range("L$r").Value = -0.1343832020997375
    $ws.Range("U$r").Value = 1.26
    $ws.Range("V$r").Value = 0.08513513513513513
    $ws.Range("W$r").Value = -0.2828729281767955
    $ws.Range("X$r").Value = 0.2482054432321119
    $ws.Range("Y$r").Value = -0.5310783714089075
    $ws.Range("Z$r").Value = 1.385908115383216
    $ws.Range("AA$r").Value = -0.107671601615074
    $ws.Range("AB$r").Value = 0.1845895865242348
    $ws.Range("AC$r").Value = -0.2922611881393088
    $ws.Range("AD$r").Value = 10.2
    $ws.Range("AF$r").Value = 10.2
    $ws.Range("AG$r").Value = 8.94
    $ws.Range("AH$r").Value = 0.408
    $ws.Range("AI$r").Value = 0.4415584415584415
    $ws.Range("AJ$r").Value = 0.3765796124684077
    $ws.Range("AK$r").Value = 0.4093406593406593
    $ws.Range("AL$r").Value = 1.9
    $ws.Range("AM$r").Value = 1.9
    $ws.Range("AN$r").Value = -5.666666666666666
    $ws.Range("AO$r").Value = -1.557894736842105
    $ws.Range("AP$r").Value = -4.966666666666666
    $ws.Range("AQ$r").Value = -1.557894736842105
}
